$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Trends Status" sheet — updated trend-status breakdown numbers
# ---------------------------------------------------------------------------
$wsTrends = $wb.Worksheets.Item("Trends Status")

# Rapid Decline
$wsTrends.Range("B2").Value = 0
$wsTrends.Range("C2").Value = 2
$wsTrends.Range("D2").Value = 0
$wsTrends.Range("E2").Value = 66.7

# Decline
$wsTrends.Range("B3").Value = 0
$wsTrends.Range("C3").Value = 0
$wsTrends.Range("D3").Value = 0
$wsTrends.Range("E3").Value = 0

# Stable
$wsTrends.Range("B4").Value = 0
$wsTrends.Range("C4").Value = 1
$wsTrends.Range("D4").Value = 0
$wsTrends.Range("E4").Value = 33.3

# Increase
$wsTrends.Range("B5").Value = 0
$wsTrends.Range("C5").Value = 0
$wsTrends.Range("D5").Value = 0
$wsTrends.Range("E5").Value = 0

# Rapid Increase
$wsTrends.Range("B6").Value = 1
$wsTrends.Range("C6").Value = 0
$wsTrends.Range("D6").Value = 100
$wsTrends.Range("E6").Value = 0

# Trend Inconclusive
$wsTrends.Range("B7").Value = 30
$wsTrends.Range("C7").Value = 60

# Insufficient Data
$wsTrends.Range("B8").Value = 477
$wsTrends.Range("C8").Value = 445

# ---------------------------------------------------------------------------
# 2. "Priority Status" sheet — updated species counts
# ---------------------------------------------------------------------------
$wsPriority = $wb.Worksheets.Item("Priority Status")
$wsPriority.Range("B2").Value = 103
$wsPriority.Range("B3").Value = 286
$wsPriority.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# 3. "Species qualification" sheet — updated label + counts
# ---------------------------------------------------------------------------
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("A2").Value = "SoIB Assessment"
$wsQual.Range("B2").Value = 508

$wsQual.Range("B3").Value = 31
$wsQual.Range("C3").Value = 1

$wsQual.Range("B4").Value = 63
$wsQual.Range("C4").Value = 3

# ---------------------------------------------------------------------------
# 4. Rename "High Priority break-up" -> "Interannual update - High Pri"
#    and update its figures; then add a new sheet right after it,
#    "Major update - High Priority ", carrying the ORIGINAL figures that
#    used to live on the "High Priority break-up" sheet.
# ---------------------------------------------------------------------------
$wsInter = $wb.Worksheets.Item("High Priority break-up")

# Create the new sheet right after it first, copying the current (pre-edit)
# values, which are the values the old "High Priority break-up" sheet used
# to hold.
$wsMajor = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsInter)
$wsMajor.Name = "Major update - High Priority "

$wsMajor.Range("A1").Value = "Break-up"
$wsMajor.Range("B1").Value = "High Species (no.)"
$wsMajor.Range("C1").Value = "High Species (perc.)"
$wsMajor.Range("D1").Value = "New High Species (no.)"
$wsMajor.Range("E1").Value = "New High Species (perc.)"

$wsMajor.Range("A2").Value = "Trend New"
$wsMajor.Range("B2").Value = 4
$wsMajor.Range("C2").Value = 25
$wsMajor.Range("D2").Value = 4
$wsMajor.Range("E2").Value = 25

$wsMajor.Range("A3").Value = "IUCN"
$wsMajor.Range("B3").Value = 12
$wsMajor.Range("C3").Value = 75
$wsMajor.Range("D3").Value = 12
$wsMajor.Range("E3").Value = 75

# Now rename the original sheet and overwrite its figures with the new ones.
$wsInter.Name = "Interannual update - High Pri"

$wsInter.Range("B2").Value = 74
$wsInter.Range("C2").Value = 71.8
$wsInter.Range("D2").Value = 74
$wsInter.Range("E2").Value = 78.7

$wsInter.Range("B3").Value = 29
$wsInter.Range("C3").Value = 28.2
$wsInter.Range("D3").Value = 20
$wsInter.Range("E3").Value = 21.3
